$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 983 (the existing rows 983:1044 shift down to 987:1048,
# matching the dimension change from A1:R1044 to A1:R1048).
$ws.Rows("983:986").Insert()

# Populate the 4 newly inserted rows with the new "Femacal de La Calera - Cebolla" records.
# Columns: A Mercado ID, B Mercado, C Region, D Fecha (serial date, same number style as
# the surrounding rows), E Codreg, F Categoria ID, G Categoria, H Variedad, I Calidad,
# J Volumen, K Precio minimo, L Precio maximo, M Precio promedio ponderado,
# N Unidad de comercializacion, O Origen, P Precio $/Kg, Q Kg o Unidades, R Clasificacion.

$rows = @(
    @{ Row = 983; D = 44610; H = "Morada(o)";        I = "1a (cosecha)"; J = 55;  K = 10000; L = 10000; M = 10000; O = "Región de Arica y Parinacota"; P = 556 },
    @{ Row = 984; D = 44610; H = "Morada(o)";        I = "2a (cosecha)"; J = 50;  K = 7000;  L = 7000;  M = 7000;  O = "Región de Arica y Parinacota"; P = 389 },
    @{ Row = 985; D = 44610; H = "Sin especificar";  I = "1a (cosecha)"; J = 133; K = 5000;  L = 5500;  M = 5256;  O = "Provincia de Quillota";        P = 292 },
    @{ Row = 986; D = 44610; H = "Sin especificar";  I = "2a (cosecha)"; J = 60;  K = 4500;  L = 4500;  M = 4500;  O = "Provincia de Quillota";        P = 250 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value2  = 3
    $ws.Cells.Item($n, 2).Value2  = "Femacal de La Calera"
    $ws.Cells.Item($n, 3).Value2  = "Coquimbo"
    $ws.Cells.Item($n, 4).Value2  = $r.D
    $ws.Cells.Item($n, 5).Value2  = 5
    $ws.Cells.Item($n, 6).Value2  = 100112004
    $ws.Cells.Item($n, 7).Value2  = "Cebolla"
    $ws.Cells.Item($n, 8).Value2  = $r.H
    $ws.Cells.Item($n, 9).Value2  = $r.I
    $ws.Cells.Item($n, 10).Value2 = $r.J
    $ws.Cells.Item($n, 11).Value2 = $r.K
    $ws.Cells.Item($n, 12).Value2 = $r.L
    $ws.Cells.Item($n, 13).Value2 = $r.M
    $ws.Cells.Item($n, 14).Value2 = "`$/malla 18 kilos"
    $ws.Cells.Item($n, 15).Value2 = $r.O
    $ws.Cells.Item($n, 16).Value2 = $r.P
    $ws.Cells.Item($n, 17).Value2 = 18
    $ws.Cells.Item($n, 18).Value2 = "Hortaliza"
}

Write-Output "Done"
